$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121, shifting rows 121:196 down to 122:197
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new data
$ws.Cells.Item(121, 1).Value = 3
$ws.Cells.Item(121, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44438
$ws.Cells.Item(121, 5).Value = 5
$ws.Cells.Item(121, 6).Value = 100112032
$ws.Cells.Item(121, 7).Value = "Zapallo italiano"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 218
$ws.Cells.Item(121, 11).Value = 12000
$ws.Cells.Item(121, 12).Value = 13000
$ws.Cells.Item(121, 13).Value = 12518
$ws.Cells.Item(121, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 16).Value = 179
$ws.Cells.Item(121, 17).Value = 70
$ws.Cells.Item(121, 18).Value = "Hortaliza"
